{"js": "// Updated driving privileges template.\n// Change: \" below, to and from home sta\" -> \" below, to and from home, including reasonable commute time, sta\"\n// (i.e. insert \", including reasonable commute time,\" right after \"...from home\" and\n// before the \"sta\" that begins \"starting on {{ plea_trial_date }}.\")\n\nconst body = context.document.body;\n\nconst searchText = \" below, to and from home sta\";\nconst results = body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find target text: ${searchText}`);\n}\n\nconst target = results.items[0];\nconst replacement = \" below, to and from home, including reasonable commute time, sta\";\ntarget.insertText(replacement, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Updated driving privileges template.\n# Change: \" below, to and from home sta\" -> \" below, to and from home, including reasonable commute time, sta\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \" below, to and from home sta\"\n$find.Replacement.Text = \" below, to and from home, including reasonable commute time, sta\"\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n$find.Execute(\n    [ref]$find.Text,\n    [ref]$true,            # MatchCase\n    [ref]$false,           # MatchWholeWord\n    [ref]$false,           # MatchWildcards\n    [ref]$false,           # MatchSoundsLike\n    [ref]$false,           # MatchAllWordForms\n    [ref]$true,            # Forward\n    [ref]1,                # Wrap (wdFindContinue)\n    [ref]$false,           # Format\n    [ref]$find.Replacement.Text,\n    [ref]2                 # Replace (wdReplaceAll)\n) | Out-Null\n"}
